# CW3M_McKenzie.xlsx regression-testing workbook update:
# Roll back a change made prematurely in a previous revision (Flow_McKenzie.xml)
# and correct a bug in the new-cloudiness calculation (ReachRouting.cpp).
# This adds a new averaged "Baseline 2010-18 C125" summary row (row 40) that
# now averages rows 31:39 (i.e. includes the previous C123 average row, 39,
# in the window) for the new corrected run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row 40: label cells -------------------------------------------
$ws.Range("A40").Value2 = "CW3M"
$ws.Range("B40").Value2 = "Baseline 2010-18 C125"
$ws.Range("C40").Value2 = "2010-18"
$ws.Range("S40").Value2 = "2010-18"

# --- new row 40: averaged measurement columns ---------------------------
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols) {
    $ws.Range($col + "40").Formula = "=AVERAGE(" + $col + "31:" + $col + "39)"
}

# --- formatting: plain numeric styles (no highlight), same as other data rows
$ws.Range("D40:N40").NumberFormat = "0.00"
$ws.Range("O40:P40").NumberFormat = "0"
$ws.Range("Q40").NumberFormat = "0.00"
$ws.Range("R40").NumberFormat = "0.000000"

# --- freeze header row and scroll/select like the saved view ------------
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("S40").Select()
